$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.892.27'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.889.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  -0.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.89'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.74%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3127'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07187'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08542'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7653'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.894.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.368'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '93.67'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.157'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.856.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.80'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.48'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007820'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.18%  '
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9990'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.28%  '
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.148.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.011'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1648'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.430'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.036'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.460'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.537'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.499'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.098'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05451'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.246'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7429'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.000'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.702'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01954'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("E40").Value = '  -0.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4474'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.101.96'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '73.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.45%  '
$ws.Range("E44").Value = '  +1.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8528'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.80%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.687'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.871'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.78%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.001'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.034.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.74%  '
